{"js": "// Update the division problems in the practice table.\n// Each cell's \"before\" text is unique in the document, so a plain\n// body.search() for each pair unambiguously locates the single cell\n// to update; insertText(..., Word.InsertLocation.replace) swaps it in.\nconst body = context.document.body;\n\nconst pairs = [\n  [\"196\u00f75=\", \"918\u00f74=\"],\n  [\"923\u00f73=\", \"374\u00f73=\"],\n  [\"670\u00f74=\", \"221\u00f79=\"],\n  [\"115\u00f73=\", \"793\u00f76=\"],\n  [\"325\u00f79=\", \"367\u00f72=\"],\n  [\"507\u00f76=\", \"229\u00f73=\"],\n  [\"227\u00f72=\", \"600\u00f76=\"],\n  [\"887\u00f75=\", \"640\u00f76=\"],\n  [\"521\u00f76=\", \"560\u00f76=\"],\n  [\"176\u00f76=\", \"756\u00f78=\"],\n  [\"169\u00f73=\", \"359\u00f79=\"],\n  [\"725\u00f73=\", \"830\u00f73=\"],\n  [\"876\u00f79=\", \"438\u00f73=\"],\n  [\"814\u00f77=\", \"670\u00f73=\"],\n  [\"755\u00f78=\", \"904\u00f77=\"],\n  [\"325\u00f72=\", \"679\u00f75=\"],\n  [\"611\u00f75=\", \"352\u00f72=\"],\n  [\"311\u00f79=\", \"559\u00f75=\"],\n  [\"619\u00f77=\", \"635\u00f73=\"],\n  [\"613\u00f73=\", \"938\u00f78=\"],\n  [\"897\u00f73=\", \"289\u00f77=\"],\n  [\"978\u00f74=\", \"502\u00f72=\"],\n  [\"418\u00f72=\", \"822\u00f79=\"],\n  [\"778\u00f75=\", \"584\u00f73=\"],\n  [\"166\u00f73=\", \"988\u00f74=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division problems in the practice table.\n# Each cell's \"before\" text is unique in the document, so a plain\n# Find/Replace (no wildcards) for each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{old=\"196\u00f75=\"; new=\"918\u00f74=\"},\n    @{old=\"923\u00f73=\"; new=\"374\u00f73=\"},\n    @{old=\"670\u00f74=\"; new=\"221\u00f79=\"},\n    @{old=\"115\u00f73=\"; new=\"793\u00f76=\"},\n    @{old=\"325\u00f79=\"; new=\"367\u00f72=\"},\n    @{old=\"507\u00f76=\"; new=\"229\u00f73=\"},\n    @{old=\"227\u00f72=\"; new=\"600\u00f76=\"},\n    @{old=\"887\u00f75=\"; new=\"640\u00f76=\"},\n    @{old=\"521\u00f76=\"; new=\"560\u00f76=\"},\n    @{old=\"176\u00f76=\"; new=\"756\u00f78=\"},\n    @{old=\"169\u00f73=\"; new=\"359\u00f79=\"},\n    @{old=\"725\u00f73=\"; new=\"830\u00f73=\"},\n    @{old=\"876\u00f79=\"; new=\"438\u00f73=\"},\n    @{old=\"814\u00f77=\"; new=\"670\u00f73=\"},\n    @{old=\"755\u00f78=\"; new=\"904\u00f77=\"},\n    @{old=\"325\u00f72=\"; new=\"679\u00f75=\"},\n    @{old=\"611\u00f75=\"; new=\"352\u00f72=\"},\n    @{old=\"311\u00f79=\"; new=\"559\u00f75=\"},\n    @{old=\"619\u00f77=\"; new=\"635\u00f73=\"},\n    @{old=\"613\u00f73=\"; new=\"938\u00f78=\"},\n    @{old=\"897\u00f73=\"; new=\"289\u00f77=\"},\n    @{old=\"978\u00f74=\"; new=\"502\u00f72=\"},\n    @{old=\"418\u00f72=\"; new=\"822\u00f79=\"},\n    @{old=\"778\u00f75=\"; new=\"584\u00f73=\"},\n    @{old=\"166\u00f73=\"; new=\"988\u00f74=\"}\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $p.old\n    $find.Replacement.Text = $p.new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
